$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:A11").Value = "Blood-results-redacted"
$ws.Range("B3:B11").Value = "Gemma3"

$null = $ws.Range("A2:H11").Select()
